# Append 9 new "user_detail" master-data rows (ids 110021-110029) to Sheet1,
# mirroring the 20 existing rows, then update the view/selection to match
# where a user would land after typing this block of rows in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: id, uin, name, email, mobile
$newRows = @(
    @{ Row = 22; Id = 110021; Uin = 7316931025; Name = "Magdalena Weber";   Email = "magdalena.weber@xyz.com";   Mobile = 932122450 },
    @{ Row = 23; Id = 110022; Uin = 9137847236; Name = "Adrienne Hoffman";  Email = "adrienne.hoffman@xyz.com";  Mobile = 848488000 },
    @{ Row = 24; Id = 110023; Uin = 8428758532; Name = "Adrienne Mcgee";    Email = "adrienne.mcgee@xyz.com";    Mobile = 894773246 },
    @{ Row = 25; Id = 110024; Uin = 9804209494; Name = "Amare Coleman";     Email = "amare.coleman@xyz.com";     Mobile = 956554588 },
    @{ Row = 26; Id = 110025; Uin = 7105248214; Name = "Dawson Ibarra";     Email = "dawson.ibarra@xyz.com";     Mobile = 765455583 },
    @{ Row = 27; Id = 110026; Uin = 9316557128; Name = "Elvis Mcmillan";    Email = "elvis.mcmillan@xyz.com";    Mobile = 884282274 },
    @{ Row = 28; Id = 110027; Uin = 8103486949; Name = "Steve George";     Email = "steve.george@xyz.com";      Mobile = 971073663 },
    @{ Row = 29; Id = 110028; Uin = 9601932866; Name = "Colton Elliott";    Email = "colton.elliott@xyz.com";    Mobile = 809908673 },
    @{ Row = 30; Id = 110029; Uin = 9317596765; Name = "Carolyn Rodriguez"; Email = "carolyn.rodriguez@xyz.com"; Mobile = 818876429 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Id          # A: id
    $ws.Cells.Item($row, 2).Value = $r.Uin         # B: uin
    $ws.Cells.Item($row, 3).Value = $r.Name        # C: name
    $ws.Cells.Item($row, 4).Value = $r.Email       # D: email
    $ws.Cells.Item($row, 5).Value = $r.Mobile      # E: mobile
    $ws.Cells.Item($row, 6).Value = "ACT"          # F: status_code
    $ws.Cells.Item($row, 7).Value = "eng"          # G: lang_code
    $ws.Cells.Item($row, 8).Value = "PWD"          # H: last_login_method
    $c9 = $ws.Cells.Item($row, 9)
    $c9.Value = $true                              # I: is_active
    $c9.HorizontalAlignment = -4131                # xlLeft (matches existing rows' style)
    $ws.Cells.Item($row, 10).Value = "superadmin"  # J: cr_by
    $ws.Cells.Item($row, 11).Value = "now()"       # K: cr_dtimes
}

# Reposition the view like a user scrolling to / selecting the freshly typed block.
$ws.Range("A22:K30").Select()
